# Daily attendance processing - reorder "Recorded By" (column G) values so
# that "System" always appears first in the comma-separated list of
# recorders, preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -notlike "*,*") { continue }

    $parts = $value -split ",\s*"
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $idx = -1
    for ($i = 0; $i -lt $trimmed.Count; $i++) {
        if ($trimmed[$i].Equals("System")) {
            $idx = $i
            break
        }
    }

    if ($idx -gt 0) {
        $newOrder = @("System")
        for ($i = 0; $i -lt $trimmed.Count; $i++) {
            if ($i -ne $idx) { $newOrder += $trimmed[$i] }
        }
        $newValue = [string]::Join(", ", $newOrder)
        $cell.Value2 = $newValue
    }
}
